# Update cryptocurrency price/volume figures in the "cryptos" sheet.
# Values are stored as literal text (e.g. "294.83", "1.38%") in the source
# workbook, so each target cell is switched to the "@" (Text) number format
# before the new literal is written -- this stops the engine from re-parsing
# number- or percent-looking strings into numeric/percentage values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "294.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.38%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.96%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.932"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.07%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07423"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.60%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.251"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "25.82%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.748"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.15%"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.746"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.14%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9143"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.23%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09076"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "17.47%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1706"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.39%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08301"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.44%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03134"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.73%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1003"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.04%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001511"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.87%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005823"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.14%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.503"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.95%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.075"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.24%"
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.52%"
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.59%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.980"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.43%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2103"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.26%"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04547"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.73%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.42%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004601"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "14.81%"
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.22%"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003401"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-95.48%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.00%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04529"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.37%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007320"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.64%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009852"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "25.12%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1330"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001904"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.76%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009137"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006191"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.26%"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.26%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.212"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.57%"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002004"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-33.16%"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.26%"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.26%"
